{"js": "// Fix typo: \"getElementByClassName\" -> \"getElementsByClassName\"\n// (the author inserted a missing \"s\" right after \"getElement\").\nconst results = context.document.body.search(\"getElementByClassName\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"getElementByClassName not found\");\n}\n\n// There is exactly one occurrence in the document; replace its text in\n// place so the surrounding spell/grammar-check proofing marks stay put.\nresults.items[0].insertText(\"getElementsByClassName\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix typo: \"getElementByClassName\" -> \"getElementsByClassName\"\n# (a missing \"s\" was inserted right after \"getElement\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"getElementByClassName\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"getElementsByClassName\"\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
